$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = -0.6303854875283446
$ws.Range("H2").Value = -0.6303854875283446
$ws.Range("I2").Value = -1.977324263038549
$ws.Range("J2").Value = -1.909292561909252
$ws.Range("K2").Value = 14.86
$ws.Range("L2").Value = 3.369614512471656
$ws.Range("M2").Value = 39.1
$ws.Range("N2").Value = 0.03009544334975369
$ws.Range("O2").Value = 2.631224764468371
$ws.Range("P2").Value = 39.1
$ws.Range("Q2").Value = 0.03009544334975369
$ws.Range("R2").Value = 2.631224764468371
$ws.Range("U2").Value = 11.46
$ws.Range("V2").Value = 0.008820812807881773
$ws.Range("W2").Value = -0.07985464745633049
$ws.Range("X2").Value = 0.0528070435897599
$ws.Range("Y2").Value = -0.1326616910460904
$ws.Range("Z2").Value = 0.004877417556712422
$ws.Range("AA2").Value = -0.2701700525291664
$ws.Range("AB2").Value = 0.04897064212014163
$ws.Range("AC2").Value = -0.3191406946493081
$ws.Range("AD2").Value = 306.79
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 306.79
$ws.Range("AG2").Value = 295.33
$ws.Range("AH2").Value = 0.1910285867284354
$ws.Range("AI2").Value = 0.2926854864099066
$ws.Range("AJ2").Value = 0.1852144519074586
$ws.Range("AK2").Value = 0.2848668409325475
$ws.Range("AL2").Value = 5.722
$ws.Range("AM2").Value = -29.278
$ws.Range("AN2").Value = -97.70382165605096
$ws.Range("AO2").Value = -1.52394267738553
$ws.Range("AP2").Value = -94.05414012738855
$ws.Range("AQ2").Value = 0.2978345515404058
$ws.Range("I3").Value = -1.626470588235294
$ws.Range("J3").Value = -1.514550087361678
$ws.Range("K3").Value = 18.8
$ws.Range("L3").Value = 5.529411764705883
$ws.Range("M3").Value = 39.1
$ws.Range("N3").Value = 0.0370159992426394
$ws.Range("O3").Value = 2.079787234042553
$ws.Range("P3").Value = 39.1
$ws.Range("Q3").Value = 0.0370159992426394
$ws.Range("R3").Value = 2.079787234042553
$ws.Range("U3").Value = 9.52
$ws.Range("V3").Value = 0.009012591119946985
$ws.Range("W3").Value = 0.03062886933854676
$ws.Range("X3").Value = 0.05737577465592213
$ws.Range("Y3").Value = -0.02674690531737537
$ws.Range("Z3").Value = 0.003785348474727232
$ws.Range("AA3").Value = -0.005733099863092522
$ws.Range("AB3").Value = 0.04979083400179928
$ws.Range("AC3").Value = -0.05552393386489179
$ws.Range("AD3").Value = 303.1
$ws.Range("AF3").Value = 303.1
$ws.Range("AG3").Value = 293.58
$ws.Range("AH3").Value = 0.2229660144181256
$ws.Range("AI3").Value = 0.2948156794086179
$ws.Range("AJ3").Value = 0.2174859987554449
$ws.Range("AK3").Value = 0.2882247835221583
$ws.Range("AL3").Value = 4.81
$ws.Range("AM3").Value = -30.19
$ws.Range("AO3").Value = -1.14968814968815
$ws.Range("AQ3").Value = 0.1831732361709175
$ws.Range("G4").Value = -2.752475247524752
$ws.Range("H4").Value = -2.752475247524752
$ws.Range("I4").Value = -3.158415841584159
$ws.Range("J4").Value = -3.158415841584159
$ws.Range("K4").Value = -3.94
$ws.Range("L4").Value = -3.900990099009901
$ws.Range("U4").Value = 1.94
$ws.Range("V4").Value = 0.007986825854261012
$ws.Range("W4").Value = -0.1903381642512077
$ws.Range("X4").Value = 0.04823831252359767
$ws.Range("Y4").Value = -0.2385764767748054
$ws.Range("Z4").Value = 0.1692642869113457
$ws.Range("AA4").Value = -0.5346070051952403
$ws.Range("AB4").Value = 0.04897064212014163
$ws.Range("AC4").Value = -0.5827574554337243
$ws.Range("AD4").Value = 3.69
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 3.69
$ws.Range("AG4").Value = 1.75
$ws.Range("AH4").Value = 0.0149641104667667
$ws.Range("AI4").Value = 0.1836734693877551
$ws.Range("AJ4").Value = 0.00715307582260372
$ws.Range("AK4").Value = 0.09641873278236915
$ws.Range("AL4").Value = 0.912
$ws.Range("AM4").Value = 0.912
$ws.Range("AN4").Value = -1.17515923566879
$ws.Range("AO4").Value = -3.49780701754386
$ws.Range("AP4").Value = -0.5573248407643312
$ws.Range("AQ4").Value = -3.49780701754386
